# Applies the TestData/testdata.xlsx edit described by the commit:
# "Enhanced code to handle website responses while running login test
#  with valid user. Modified LoginPage.java and LoginTest.java"
#
# Net effect on the worksheet data:
#   - Row 2 (A2/B2) now holds the admin credentials (admin@yourstore.com / admin)
#     and gets an explicit 15pt row height.
#   - Row 3 (A3/B3) now holds admin@yourstore.com / pass124
#   - Row 4 (A4/B4) now holds user1@test.com / admin
#   - Row 5 (A5/B5) now holds user2@test.com / pass126
#   - Row 6 is removed entirely (dimension shrinks from A1:B6 to A1:B5)
#   - Stale per-cell hyperlinks are replaced with just two: A4 -> user1@test.com,
#     A5 -> user2@test.com (mailto: links), matching the surviving cell content.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the surviving data rows -------------------------------------
$ws.Range("A2").Value = "admin@yourstore.com"
$ws.Range("B2").Value = "admin"

$ws.Range("A3").Value = "admin@yourstore.com"
$ws.Range("B3").Value = "pass124"

$ws.Range("A4").Value = "user1@test.com"
$ws.Range("B4").Value = "admin"

$ws.Range("A5").Value = "user2@test.com"
$ws.Range("B5").Value = "pass126"

# --- Drop the now-unused last row (admin / admin row 6) -----------------
$ws.Rows.Item(6).Delete()

# --- Row 2 gets an explicit row height in the edited workbook -----------
$ws.Rows.Item(2).RowHeight = 15

# --- Rebuild hyperlinks: only A4 and A5 keep a mailto: link now ---------
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A4"), "mailto:user1@test.com")
$ws.Hyperlinks.Add($ws.Range("A5"), "mailto:user2@test.com")

# Hyperlinks.Add() re-stamps a fresh cell style; put column A's hyperlink
# cells back on the shared "Hyperlink" style used throughout the sheet.
$ws.Range("A4").Style = "Hyperlink"
$ws.Range("A5").Style = "Hyperlink"

Write-Host "testdata.xlsx updated: admin login rows consolidated, stale rows/hyperlinks removed"
